# Fix esm bin bug: shrink the test fixture down to the actual 3x4 data block
# that is used by the "empty data row" test case. The worksheet originally
# padded its real A1:C4 table out to A1:G21 with empty, styled filler
# cells/columns - those are no longer needed, so remove them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the empty filler rows 5:21 (shifts nothing below them, they are
# the last rows of data), leaving only the 4 real data rows.
$ws.Range("A5:G21").EntireRow.Delete()

# Remove the empty filler columns D:G, leaving only the 3 real data columns.
$ws.Range("D1:G4").EntireColumn.Delete()

Write-Output ("New used range: " + $ws.UsedRange.Address())
